# CHP and heat_pump linearized
# Remove the solar_th, pvt and gas_boiler columns, keeping only
# net1, CHP1 and heat_pump1. The heat_pump1 column (previously F)
# moves into column D, and columns E, F, G are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move heat_pump1 values (currently in column F) into column D
$ws.Range("D1").Value = $ws.Range("F1").Value()
$ws.Range("D2").Value = $ws.Range("F2").Value()
$ws.Range("D3").Value = $ws.Range("F3").Value()

# Clear the now-unused columns E, F and G entirely
$ws.Range("E1:G3").Clear()
